$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental value -> false (leading apostrophe forces literal text storage
# instead of Excel auto-converting "false" to a Boolean). Re-apply the
# original cell's formatting afterwards so the quote-prefix marker doesn't
# change the cell's effective style.
$ws.Range("B7").Value = "'false"
$ws.Range("B18").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null

# Date value updated
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Description value filled in
$ws.Range("B17").Value = "Codes for trend direction indicators"
